$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 31, shifting rows 31..49 down to 32..50
$ws.Rows.Item(31).Insert()

# Populate the new row with the "Eqlimy Gate" entry
$ws.Range("A31").Value = "Eqlimy Gate"
$ws.Range("B31").Value = "169.254.1.81"

# Match the resulting selection from the edit
$ws.Range("B32").Select()
